$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E86").Value = 3.680259274294004
$ws.Range("F86").Value = 5.271211974840319
$ws.Range("H86").Value = 0.08878512403546765
$ws.Range("I86").Value = 3.846976579526866
$ws.Range("E87").Value = 2.347689841918489
$ws.Range("F87").Value = 5.843826550722718
$ws.Range("H87").Value = 0.1440987094069509
$ws.Range("I87").Value = 2.615974823067429
$ws.Range("E88").Value = -8.976776321850711
$ws.Range("F88").Value = 0.5682929652578221
$ws.Range("H88").Value = 0.02453766901004156
$ws.Range("I88").Value = -8.560149267884512
$ws.Range("E89").Value = -3.978309616519428
$ws.Range("F89").Value = -1.731784205539412
$ws.Range("H89").Value = 0.1125882280211566
$ws.Range("I89").Value = -3.702861693241494
$ws.Range("E90").Value = -2.102616373543336
$ws.Range("F90").Value = -3.177503117498747
$ws.Range("H90").Value = -0.2552882542189747
$ws.Range("I90").Value = -1.981230449483548
$ws.Range("E91").Value = -0.731437487636884
$ws.Range("F91").Value = -3.94728494988759
$ws.Range("H91").Value = -0.2052199120795587
$ws.Range("I91").Value = -0.6746413364222819
$ws.Range("E92").Value = -0.1989929386189175
$ws.Range("F92").Value = -1.752839104079642
$ws.Range("H92").Value = -0.2006453597222053
$ws.Range("I92").Value = -0.1416849346244687
$ws.Range("E93").Value = 0.5533475792865887
$ws.Range("F93").Value = -0.6199248051281374
$ws.Range("H93").Value = -0.2672854350458368
$ws.Range("I93").Value = 0.6485735476804869
$ws.Range("E94").Value = 0.1847359823484989
$ws.Range("F94").Value = -0.04808671615517868
$ws.Range("H94").Value = -0.3076969134591389
$ws.Range("I94").Value = 0.3653291761176868
$ws.Range("E95").Value = -0.02586288097359801
$ws.Range("F95").Value = 0.1283069355106428
$ws.Range("H95").Value = -0.3349863823792135
$ws.Range("I95").Value = 0.1455320217330877
$ws.Range("E96").Value = -0.6911101911302666
$ws.Range("F96").Value = 0.00527762238280554
$ws.Range("H96").Value = -0.2875140414930297
$ws.Range("I96").Value = -0.5286778324510912
$ws.Range("E97").Value = -0.7109463283291764
$ws.Range("F97").Value = -0.3107958545211357
$ws.Range("H97").Value = -0.308388168623663
$ws.Range("I97").Value = -0.5556860244073824
$ws.Range("E98").Value = -0.6385110679305853
$ws.Range("F98").Value = -0.5166076170909067
$ws.Range("H98").Value = -0.1242078811297451
$ws.Range("I98").Value = -0.478814596627843
$ws.Range("E99").Value = -0.5636118526484968
$ws.Range("F99").Value = -0.6510448600096315
$ws.Range("H99").Value = -0.1131410748646612
$ws.Range("I99").Value = -0.4139735871137652
$ws.Range("E100").Value = -0.5055187398544909
$ws.Range("F100").Value = -0.6046469971906876
$ws.Range("H100").Value = -0.0987316757413821
$ws.Range("I100").Value = -0.3752019208453369
$ws.Range("E101").Value = -0.4631739489024822
$ws.Range("F101").Value = -0.5427039023340141
$ws.Range("H101").Value = -0.08063665727845436
$ws.Range("I101").Value = -0.3481855133441815
$ws.Range("E102").Value = -0.2673038422047105
$ws.Range("F102").Value = -0.4499020959025454
$ws.Range("H102").Value = -0.05723682191681605
$ws.Range("I102").Value = -0.152780448922243
$ws.Range("E103").Value = -0.08528450216597984
$ws.Range("F103").Value = -0.3303202582819161
$ws.Range("H103").Value = -0.07587101817681705
$ws.Range("I103").Value = 0.0156357572484178
$ws.Range("E104").Value = 0.05812151651325217
$ws.Range("F104").Value = -0.1894101941899803
$ws.Range("H104").Value = -0.05486731123815808
$ws.Range("I104").Value = 0.1722365894389535
$ws.Range("E105").Value = 0.2047767174390268
$ws.Range("F105").Value = -0.02242252760460312
$ws.Range("H105").Value = -0.04526669706578193
$ws.Range("I105").Value = 0.309879854442888
$ws.Range("E106").Value = 0.0722737176454564
$ws.Range("F106").Value = 0.06247186235793861
$ws.Range("H106").Value = -0.03475888848409577
$ws.Range("I106").Value = 0.1703023180344354
$ws.Range("E107").Value = -0.08670008486227136
$ws.Range("F107").Value = 0.06211796668386572
$ws.Range("H107").Value = -0.05443444842384374
$ws.Range("I107").Value = -0.0002624824500492945
$ws.Range("E108").Value = -0.2164226799810804
$ws.Range("F108").Value = -0.006518082439717431
$ws.Range("H108").Value = -0.02262340175337374
$ws.Range("I108").Value = -0.1309253276396118
$ws.Range("E109").Value = -0.2642465696146233
$ws.Range("F109").Value = -0.1237739042031299
$ws.Range("H109").Value = -0.03324402107510648
$ws.Range("I109").Value = -0.2000054108628525
$ws.Range("E110").Value = -0.2873397052640327
$ws.Range("F110").Value = -0.2136772599305022
$ws.Range("H110").Value = -0.0168875205374413
$ws.Range("I110").Value = -0.2028893577346539
$ws.Range("E111").Value = -0.2968015602886593
$ws.Range("F111").Value = -0.2662026287870992
$ws.Range("H111").Value = -0.02510931265696181
$ws.Range("I111").Value = -0.2046584397615549
$ws.Range("E112").Value = -0.2949165131893148
$ws.Range("F112").Value = -0.2858260870891578
$ws.Range("H112").Value = -0.03521861439247089
$ws.Range("I112").Value = -0.2245933613779398
$ws.Range("E113").Value = -0.3529101000271226
$ws.Range("F113").Value = -0.3079919696922827
$ws.Range("H113").Value = -0.03198795660782582
$ws.Range("I113").Value = -0.2557517862576807
$ws.Range("E114").Value = -0.2214324879033238
$ws.Range("F114").Value = -0.2915151653521054
$ws.Range("H114").Value = -0.04142977303710259
$ws.Range("I114").Value = -0.1432161542508758
$ws.Range("E115").Value = -0.0634273571989695
$ws.Range("F115").Value = -0.2331716145796829
$ws.Range("H115").Value = -0.02673567703335961
$ws.Range("I115").Value = 0.03259648744162259
$ws.Range("E116").Value = 0.1107574097593994
$ws.Range("F116").Value = -0.1317531338425044
$ws.Range("H116").Value = -0.04797402840468178
$ws.Range("I116").Value = 0.1953754212427777
$ws.Range("E117").Value = 0.2265775790995693
$ws.Range("F117").Value = 0.01311878593916858
$ws.Range("H117").Value = -0.05239837697603862
$ws.Range("I117").Value = 0.3157760010203308
$ws.Range("E118").Value = -0.1232732898647901
$ws.Range("F118").Value = 0.037658585448802
$ws.Range("H118").Value = -0.04105292735169145
$ws.Range("I118").Value = -0.0116742036345335
$ws.Range("E119").Value = -0.3846861531926303
$ws.Range("F119").Value = -0.0426561135496132
$ws.Range("H119").Value = -0.05583481397406959
$ws.Range("I119").Value = -0.2871788880414762
$ws.Range("E120").Value = -0.6223825479344072
$ws.Range("F120").Value = -0.2259411029730649
$ws.Range("H120").Value = -0.04848284161678198
$ws.Range("I120").Value = -0.502856511207016
$ws.Range("E121").Value = -0.7702716412667178
$ws.Range("F121").Value = -0.4751534080646366
$ws.Range("H121").Value = -0.06656174414540586
$ws.Range("I121").Value = -0.6626803170988551
$ws.Range("E122").Value = -0.3313267991198693
$ws.Range("F122").Value = -0.5271667853784064
$ws.Range("H122").Value = -0.05873691752142227
$ws.Range("I122").Value = -0.220803444024164
$ws.Range("E123").Value = -0.003709895504907329
$ws.Range("F123").Value = -0.4319227209564757
$ws.Range("H123").Value = -0.04850694202569691
$ws.Range("I123").Value = 0.1235952038025795
$ws.Range("E124").Value = 0.2574326603467768
$ws.Range("F124").Value = -0.2119689188861796
$ws.Range("H124").Value = -0.06831112276902784
$ws.Range("I124").Value = 0.3767978866662327
$ws.Range("E125").Value = 0.4366669931689149
$ws.Range("F125").Value = 0.08976573972272853
$ws.Range("H125").Value = -0.06930120098861753
$ws.Range("I125").Value = 0.5541311136455156
$ws.Range("F126").Value = 0.2817641877949246
$ws.Range("I126").Value = 0.437669051895288
